$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing cell C2: 3.2 -> 4.2
$ws.Range("C2").Value = 4.2

# Append new ranking rows 477-495
$newRows = @(
  @(477, "אביב ואסקז", 1),
  @(478, "ליהי בראל", 1),
  @(479, "איתי הראל", 1),
  @(480, "שלו דיין", 1),
  @(481, "אורי שטרנברג", 1),
  @(482, "שלו דיין", 6),
  @(483, "אורי שטרנברג", 6),
  @(484, "הגר אגמון", 1),
  @(485, "דן פימה", 1),
  @(486, "איתי הראל", 1),
  @(487, "הילה שולויס", 1),
  @(488, "איתי בסטקר", 1),
  @(489, "קרן רינת פביאן", 1),
  @(490, "יולי קזמה", 1),
  @(491, "יהלי דוייב", 1),
  @(492, "אן מרש", 1),
  @(493, "תאיו ורד", 1),
  @(494, "יהלי דוייב", 6),
  @(495, "הגר אגמון", 6)
)

foreach ($r in $newRows) {
  $rowNum = $r[0]
  $name = $r[1]
  $pts = $r[2]
  $ws.Cells.Item($rowNum, 1).Value = $name
  $ws.Cells.Item($rowNum, 2).Value = $pts
}

# Update selection / scroll position to match the saved view state
$ws.Range("D5").Select() | Out-Null

